# "Add files via upload" — replaces the two-sheet COVID tracker (SVTeam +
# Examples) with the single consolidated sheet used going forward:
#   - drop the stale "Examples" scratch sheet
#   - rename "SVTeam" -> "Sheet1"
#   - the sheet no longer needs the header AutoFilter
#   - leave the cursor where the author last left it (B22)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets("SVTeam")
$examples = $wb.Worksheets("Examples")
$examples.Delete() | Out-Null

$ws.Name = "Sheet1"

# The autofilter dropdowns on row 1 are removed.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Row 1's header height settles back to an (non-custom) auto height.
$ws.Rows.Item(1).RowHeight = 48.6

# Restore the last-used selection.
$ws.Range("B22").Select() | Out-Null
